$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.661.77"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "3.551.95"
$ws.Range("E3").Value = "  -1.47%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "197.13"
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("D6").Value = "585.77"
$ws.Range("E6").Value = "  -3.16%  "

$ws.Range("D7").Value = "0.612"
$ws.Range("E7").Value = "  -2.21%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "0.207"
$ws.Range("E9").Value = "  +0.69%  "

$ws.Range("D10").Value = "0.632"
$ws.Range("E10").Value = "  -2.49%  "

$ws.Range("D11").Value = "'52.10"
$ws.Range("E11").Value = "  -3.15%  "

$ws.Range("D12").Value = "'0.0000288"
$ws.Range("E12").Value = "  -4.73%  "

$ws.Range("D13").Value = "9.26"
$ws.Range("E13").Value = "  -3.41%  "

$ws.Range("D14").Value = "4.109.52"
$ws.Range("E14").Value = "  -1.61%  "

$ws.Range("D15").Value = "671.31"
$ws.Range("E15").Value = "  +12.69%  "

$ws.Range("D16").Value = "69.646.96"
$ws.Range("E16").Value = "  -0.96%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.549.08"
$ws.Range("E17").Value = "  -1.96%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "12.53"
$ws.Range("E18").Value = "  -4.37%  "

$ws.Range("E19").Value = "  -0.71%  "

$ws.Range("D20").Value = "18.48"
$ws.Range("E20").Value = "  -3.24%  "

$ws.Range("D21").Value = "0.966"
$ws.Range("E21").Value = "  -3.06%  "

$ws.Range("D22").Value = "18.06"
$ws.Range("E22").Value = "  +1.67%  "

$ws.Range("D23").Value = "5.33"
$ws.Range("E23").Value = "  +2.46%  "

$ws.Range("D24").Value = "105.53"
$ws.Range("E24").Value = "  +3.40%  "

$ws.Range("D25").Value = "4.39"
$ws.Range("E25").Value = "  -4.80%  "

$ws.Range("D26").Value = "2.93"
$ws.Range("E26").Value = "  -2.85%  "

$ws.Range("D27").Value = "10.19"
$ws.Range("E27").Value = "  -5.03%  "

$ws.Range("D28").Value = "9.63"
$ws.Range("E28").Value = "  +0.17%  "

$ws.Range("D29").Value = "33.53"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("D30").Value = "4.38"
$ws.Range("E30").Value = "  -8.08%  "

$ws.Range("D31").Value = "6.77"
$ws.Range("E31").Value = "  -5.28%  "

$ws.Range("D32").Value = "11.76"
$ws.Range("E32").Value = "  -4.18%  "

$ws.Range("E33").Value = "  -4.46%  "

$ws.Range("D34").Value = "62.06"
$ws.Range("E34").Value = "  -1.99%  "

$ws.Range("D35").Value = "3.782.45"
$ws.Range("E35").Value = "  -3.16%  "

$ws.Range("D36").Value = "0.0₃0818"
$ws.Range("E36").Value = "  -8.18%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "3.74"
$ws.Range("E37").Value = "  +5.66%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").Value = "502.77"
$ws.Range("E39").Value = "  -2.94%  "

$ws.Range("D40").Value = "2.91"
$ws.Range("E40").Value = "  -5.83%  "

$ws.Range("E41").Value = "  -4.64%  "

$ws.Range("D42").Value = "0.135"
$ws.Range("E42").Value = "  +1.30%  "

$ws.Range("D43").Value = "'34.70"
$ws.Range("E43").Value = "  -6.03%  "

$ws.Range("D44").Value = "'0.0450"
$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "2.88"
$ws.Range("E45").Value = "  +0.83%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "3.37"
$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("D47").Value = "0.137"
$ws.Range("E47").Value = "  -2.48%  "

$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("D49").Value = "8.35"
$ws.Range("E49").Value = "  -3.12%  "

$ws.Range("D50").Value = "1.78"
$ws.Range("E50").Value = "  +20.55%  "

$ws.Range("D51").Value = "2.74"
$ws.Range("E51").Value = "  +64.19%  "
